$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price column (D) values look numeric (e.g. "1.0000", "0.06030") but must
# remain literal text to preserve exact formatting/trailing zeros, matching
# the source scrape. Mark each one as Text before writing so Excel does not
# silently coerce the string into a number (NumberFormat on a multi-area
# Range only binds the first area, so loop cell-by-cell).
$priceCells = @("D2", "D3", "D5", "D6", "D7", "D8", "D9", "D10", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D29", "D30", "D31", "D32", "D34", "D35", "D37", "D38", "D39", "D40", "D42", "D45", "D46", "D49", "D50", "D51")
foreach ($cellAddr in $priceCells) {
    $ws.Range($cellAddr).NumberFormat = "@"
}

$ws.Range("D2").Value = "27.611.51"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "1.859.96"
$ws.Range("E3").Value = "  -0.99%  "
$ws.Range("E4").Value = "  +0.84%  "
$ws.Range("D5").Value = "333.92"
$ws.Range("E5").Value = "  +0.38%  "
$ws.Range("D6").Value = "1.011"
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("D7").Value = "0.4673"
$ws.Range("E7").Value = "  -0.93%  "
$ws.Range("D8").Value = "0.3895"
$ws.Range("E8").Value = "  -1.61%  "
$ws.Range("D9").Value = "45.28"
$ws.Range("E9").Value = "  -5.31%  "
$ws.Range("D10").Value = "0.07972"
$ws.Range("E10").Value = "  -0.99%  "
$ws.Range("D11").Value = "1.0000"
$ws.Range("D12").Value = "21.63"
$ws.Range("E12").Value = "  -2.41%  "
$ws.Range("D13").Value = "1.861.25"
$ws.Range("E13").Value = "  -0.82%  "
$ws.Range("D14").Value = "5.972"
$ws.Range("E14").Value = "  -0.10%  "
$ws.Range("D15").Value = "7.221"
$ws.Range("E15").Value = "  +1.37%  "
$ws.Range("D16").Value = "1.013"
$ws.Range("E16").Value = "  +0.63%  "
$ws.Range("D17").Value = "87.94"
$ws.Range("E17").Value = "  +0.93%  "
$ws.Range("D18").Value = "0.06698"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("E19").Value = "  -0.89%  "
$ws.Range("E20").Value = "  -2.26%  "
$ws.Range("D21").Value = "1.010"
$ws.Range("E21").Value = "  +0.63%  "
$ws.Range("D22").Value = "27.588.79"
$ws.Range("D23").Value = "5.434"
$ws.Range("E23").Value = "  -1.70%  "
$ws.Range("D24").Value = "10.83"
$ws.Range("E24").Value = "  -1.65%  "
$ws.Range("D25").Value = "2.306"
$ws.Range("E25").Value = "  +0.06%  "
$ws.Range("D26").Value = "2.080.51"
$ws.Range("E26").Value = "  -0.76%  "
$ws.Range("D27").Value = "158.56"
$ws.Range("E27").Value = "  -0.46%  "
$ws.Range("E28").Value = "  -2.22%  "
$ws.Range("D29").Value = "2.123"
$ws.Range("E29").Value = "  +0.92%  "
$ws.Range("D30").Value = "5.385"
$ws.Range("E30").Value = "  -3.48%  "
$ws.Range("D31").Value = "121.14"
$ws.Range("E31").Value = "  -0.53%  "
$ws.Range("D32").Value = "0.9711"
$ws.Range("E32").Value = "  -1.22%  "
$ws.Range("D34").Value = "3.643"
$ws.Range("E34").Value = "  +1.27%  "
$ws.Range("D35").Value = "5.279"
$ws.Range("E35").Value = "  -1.46%  "
$ws.Range("E36").Value = "  -8.41%  "
$ws.Range("D37").Value = "0.06030"
$ws.Range("D38").Value = "0.02218"
$ws.Range("E38").Value = "  -1.98%  "
$ws.Range("D39").Value = "1.193"
$ws.Range("E39").Value = "  -3.09%  "
$ws.Range("D40").Value = "8.202"
$ws.Range("E40").Value = "  +0.87%  "
$ws.Range("E41").Value = "  +0.68%  "
$ws.Range("D42").Value = "0.5906"
$ws.Range("E42").Value = "  -1.89%  "
$ws.Range("E43").Value = "  -1.45%  "
$ws.Range("E44").Value = "  -1.42%  "
$ws.Range("D45").Value = "1.249"
$ws.Range("E45").Value = "  -0.45%  "
$ws.Range("D46").Value = "0.5610"
$ws.Range("E46").Value = "  -2.13%  "
$ws.Range("E47").Value = "  -0.94%  "
$ws.Range("E48").Value = "  -1.91%  "
$ws.Range("D49").Value = "3.267"
$ws.Range("E49").Value = "  -3.54%  "
$ws.Range("D50").Value = "0.06766"
$ws.Range("E50").Value = "  -2.17%  "
$ws.Range("D51").Value = "112.47"
$ws.Range("E51").Value = "  -1.53%  "
